$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rewrite the header row (old rows 1+2 collapse into a single new row 1)
# ---------------------------------------------------------------------------
$headers = @("idx", "idx2", "Name", "Date Start", "Date End", "(m3/s)", "(MW1)", "(MW2)", "(GWh) Winter", "(GWh) Summer", "(GWh) Year")

# Columns A:E (idx, idx2, Name, Date Start, Date End) keep the default/general
# formatting (Arial 10pt, General) - reset in case any of them previously
# carried another style (e.g. old E1 used the Arial-9pt header style).
$ws.Range("A1:E1").Font.Name = "Arial"
$ws.Range("A1:E1").Font.Size = 10
$ws.Range("A1:E1").NumberFormat = "General"

for ($j = 0; $j -lt $headers.Length; $j++) {
    $cell = $ws.Cells.Item(1, 1 + $j)
    $cell.Value = $headers[$j]
}

# Columns F:K (the units row) get the Arial 9pt font used
# elsewhere in the sheet, with a plain "General" number format.
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9
$ws.Range("F1:K1").NumberFormat = "General"

# ---------------------------------------------------------------------------
# 2. Shift the 34 data rows up by one (old rows 3-36 -> new rows 2-35)
# ---------------------------------------------------------------------------
$data = @(
    ,@(1, 401225, "Schwanden (F.Blumer)", 1828, 1999, 21.5, 0.9, 0.85, 1.5, 2.8, 4.3)
    ,@(2, 404000, "Ziegelbrücke", 1837, 2011, 10, 0.58, 0.58, 1.1, 2, 3.1)
    ,@(3, 401200, "Haslen", 1848, 2012, 10.56, 0.92, 0.88, 2.76, 3.24, 6)
    ,@(4, 401000, "Hätzingen", 1850, 2007, 8, 1.05, 0.88, 2.3, 4.4, 6.7)
    ,@(5, 400700, "Linthal (Spinnerei Linthal)", 1877, 2012, 15, 1.6, 1.6, 3.64, 4.5, 8.14)
    ,@(6, 401600, "Schwanden (Gemeinde)", 1899, 2005, 0.56, 1.8, 1.8, 4.5, 6, 10.5)
    ,@(7, 400600, "Linthkraft (Stiftung)", 1900, 2004, 5, 0.5, 0.4, 1.17, 1.1, 2.27)
    ,@(8, 401300, "Krauchbach", 1902, 1999, 0.8, 0.64, 0.6, 0.8, 2.3, 3.1)
    ,@(9, 403900, "Bleiche", 1903, 1980, 0.65, 2.22, 2.1, 2.2, 4.6, 6.8)
    ,@(10, 401900, "Am Löntsch", 1908, 1975, 20, 89.1, 60, 47.5, 71.6, 119.1)
    ,@(11, 401700, "Schwanden (SNE)", 1908, 2000, 5, 0.5, 0.4, 0.8, 1.5, 2.3)
    ,@(12, 401800, "Elggiskraft, Netstal", 1912, 1983, 20, 0.8, 0.8, 1.66, 2.69, 4.35)
    ,@(13, 401775, "Mühlefuhr", 1920, 2010, 22, 0.82, 0.82, 1.65, 3.85, 5.5)
    ,@(14, 400800, "Brummbach, Linthal", 1921, 2010, 0.7, 3.9, 2.44, 5, 7.6, 12.6)
    ,@(15, 402000, "Linthkraft, Netstal", 1923, 1960, 35, 1.42, 1.4, 2.5, 3.7, 6.2)
    ,@(16, 403800, "Risi", 1925, 2011, 1.4, 4.5, 4.2, 6.6, 12.1, 18.7)
    ,@(17, 401150, "Leuggelbach", 1930, 1993, 0.2, 1, 1, 2.5, 3.9, 6.4)
    ,@(18, 401400, "Schwanden (Sernf)", 1931, 2003, 8.8, 15.7, 16.5, 26, 59, 85)
    ,@(19, 401500, "Schwanden (Niederenbach)", 1931, 2004, 3.6, 28.2, 33.5, 13, 27, 40)
    ,@(20, 401100, "Luchsingen", 1940, 1976, 0.84, 3.5, 3.3, 4, 9.5, 13.5)
    ,@(21, 400950, "Diesbach (Linth)", 1945, 1997, 12, 0.93, 0.88, 2, 3.1, 5.1)
    ,@(22, 400500, "Fätschbach", 1950, $null, 3, 15, 14.5, 20.1, 54.3, 74.4)
    ,@(23, 403600, "Rütiberg (Obersee)", 1961, 2007, 0.8, 0.89, 0.93, 1.5, 3.2, 4.7)
    ,@(24, 400400, "Linthal (Limmern)", 1964, $null, 32, 34.4, 33.5, 30.7, 50.4, 81.1)
    ,@(25, 400200, "Tierfehd (Limmern)", 1964, $null, 30, 261, 255, 172, 111.7, 283.7)
    ,@(26, 400300, "Tierfehd (Hintersand)", 1964, 1968, 10, 42, 40, 11.5, 71.5, 83)
    ,@(27, 400100, "Muttsee", 1965, $null, 0.9, 4.4, 4, 6.4, 1.4, 7.8)
    ,@(28, 403700, "Rütiberg (Brändbach)", 1977, $null, 0.38, 0.85, 0.8, 1.3, 2.5, 3.8)
    ,@(29, 401250, "Empächli", 1995, $null, 0.1, 0.33, 0.32, 0.3, 0.9, 1.2)
    ,@(30, 401350, "Engi, Vorderdorf (Mühlebach)", 2009, $null, 1.1, 3.75, 3.75, 3.45, 12.65, 16.1)
    ,@(31, 400250, "Tierfehd (Umwälzwerk)", 2010, $null, 15.3, 140, 140, $null, $null, $null)
    ,@(32, 401790, "Walzmühle 2", 2011, $null, 10, 0.32, 0.32, 1.07, 1.31, 2.38)
    ,@(33, 401325, "Engi, Hinterdorf (Sernf)", 2012, $null, 7, 0.53, 0.53, 0.4, 2.1, 2.5)
    ,@(34, 400850, "Brummbach, Braunwald", 2014, $null, 0.7, 0.93, 1, 0.4, 2.8, 3.2)
)


for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i
    $row = $data[$i]

    $rowRange = $ws.Range("A" + $r + ":K" + $r)
    $rowRange.Font.Name = "Arial"
    $rowRange.Font.Size = 9

    $ws.Range("A" + $r).NumberFormat = "0"
    $ws.Range("B" + $r).NumberFormat = "0"
    $ws.Range("C" + $r).NumberFormat = "General"
    $ws.Range("D" + $r).NumberFormat = "0"
    $ws.Range("E" + $r).NumberFormat = "0"
    $ws.Range("F" + $r + ":K" + $r).NumberFormat = "0.00"

    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, 1 + $j).Value = $row[$j]
    }
}

# ---------------------------------------------------------------------------
# 3. Clear the old 34th data row (now row 36): keep it present but empty,
#    matching the blank placeholder rows below it.
# ---------------------------------------------------------------------------
$ws.Range("A36:K36").ClearContents()
$ws.Range("A36:K36").Font.Name = "Arial"
$ws.Range("A36:K36").Font.Size = 9
$ws.Range("A36:K36").NumberFormat = "0"

# ---------------------------------------------------------------------------
# 4. Update the active selection to match the new first data row.
# ---------------------------------------------------------------------------
$ws.Range("A2:K2").Select()
